$d = $word.ActiveDocument

# 1) Fix spacing in first line "|: : :" -> "| : :"
$d.Content.Find.Execute("|: : :", $true, $false, $false, $false, $false, $true, 1, $false, "| : :", 2)

# 2) Typo fixes in question 20 line
$d.Content.Find.Execute("tied fo four ropes", $true, $false, $false, $false, $false, $true, 1, $false, "tied to four ropes", 2)
$d.Content.Find.Execute("different materials J", $true, $false, $false, $false, $false, $true, 1, $false, "different materiafs J", 2)

# 3) Remove the paragraph that contains the inline image (drawing)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Delete()
        break
    }
}

# 4) Replace "forewing" phrase with garbled OCR text including a left double quote char
$replacement4 = "Which of the " + [char]0x201C + "nnn statements is definitely true?"
$d.Content.Find.Execute("Which of the forewing statements is definitely true?", $true, $false, $false, $false, $false, $true, 1, $false, $replacement4, 2)

# 5) Fix option (1) line
$d.Content.Find.Execute("(1) Aand Bonly .", $true, $false, $false, $false, $false, $true, 1, $false, "(1}. AandBonly .", 2)

# 6) Fix option (4) line
$d.Content.Find.Execute("(4) A.B, Cand\D", $true, $false, $false, $false, $false, $true, 1, $false, "4) A.B, Cand|D", 2)

Write-Host "done"
